$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cStyle = $ws.Range("C54").Style()
Write-Host "C54 style: $cStyle"
$cVal = $ws.Range("C54").Value()
Write-Host "C54 value: $cVal"
